# Fruta / hortaliza, semanal
# Insert a new weekly record as row 151 in the data table, shifting all
# subsequent rows down by one (old row 151 -> 152, ..., old row 162 -> 163).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 151 (pushes existing rows 151-162 down to 152-163)
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row with the new record's data
$ws.Cells.Item(151, 1).Value = 9
$ws.Cells.Item(151, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(151, 3).Value = "Metropolitana"
$ws.Cells.Item(151, 4).Value = 44516
$ws.Cells.Item(151, 5).Value = 13
$ws.Cells.Item(151, 6).Value = 100112030
$ws.Cells.Item(151, 7).Value = "Poroto granado"
$ws.Cells.Item(151, 8).Value = "Sin especificar"
$ws.Cells.Item(151, 9).Value = "Primera"
$ws.Cells.Item(151, 10).Value = 25
$ws.Cells.Item(151, 11).Value = 34000
$ws.Cells.Item(151, 12).Value = 36000
$ws.Cells.Item(151, 13).Value = 34960
$ws.Cells.Item(151, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(151, 15).Value = "Perú"
$ws.Cells.Item(151, 16).Value = 1398
$ws.Cells.Item(151, 17).Value = 25
$ws.Cells.Item(151, 18).Value = "Hortaliza"
